# Auto-applied update of Leve profit-calculation columns (H:N) per scheduled Hyperion market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 732.1667
$ws.Range("I2").Value = 484.42856
$ws.Range("K2").Value = 484.42856
$ws.Range("M2").Value = -371.42856
$ws.Range("H33").Value = 569.82355
$ws.Range("I33").Value = 718.4545000000001
$ws.Range("J33").Value = 297.33334
$ws.Range("K33").Value = 718.4545000000001
$ws.Range("L33").Value = 297.33334
$ws.Range("M33").Value = -489.4545000000001
$ws.Range("N33").Value = -755.33334
$ws.Range("H43").Value = 3761.0715
$ws.Range("J43").Value = 5218.75
$ws.Range("L43").Value = 5218.75
$ws.Range("N43").Value = -5356.75
$ws.Range("H51").Value = 4631.8335
$ws.Range("I51").Value = 1800
$ws.Range("K51").Value = 1800
$ws.Range("M51").Value = -1316
$ws.Range("H116").Value = 5866.88
$ws.Range("I116").Value = 4459.077
$ws.Range("J116").Value = 7392
$ws.Range("K116").Value = 4459.077
$ws.Range("L116").Value = 7392
$ws.Range("M116").Value = -1017.077
$ws.Range("N116").Value = -14276
$ws.Range("H125").Value = 8549778
$ws.Range("I125").Value = 1404.7858
$ws.Range("K125").Value = 12643.0722
$ws.Range("M125").Value = -10183.0722
$ws.Range("H132").Value = 27028914
$ws.Range("I132").Value = 27779702
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 83339106
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -83336576
$ws.Range("N132").Value = -6560
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 48718.05
$ws.Range("I137").Value = 55515.605
$ws.Range("J137").Value = 3854.2
$ws.Range("K137").Value = 166546.815
$ws.Range("L137").Value = 11562.6
$ws.Range("M137").Value = -163996.815
$ws.Range("N137").Value = -16662.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 10249.75
$ws.Range("I3").Value = 999
$ws.Range("J3").Value = 13333.333
$ws.Range("K3").Value = 999
$ws.Range("L3").Value = 13333.333
$ws.Range("M3").Value = -884
$ws.Range("N3").Value = -13563.333
$ws.Range("H32").Value = 10072.381
$ws.Range("I32").Value = 5785.3477
$ws.Range("K32").Value = 5785.3477
$ws.Range("M32").Value = -5498.3477
$ws.Range("H41").Value = 1903.6
$ws.Range("I41").Value = 3000
$ws.Range("J41").Value = 1629.5
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 1629.5
$ws.Range("M41").Value = -2586
$ws.Range("N41").Value = -2457.5
$ws.Range("H51").Value = 33100
$ws.Range("J51").Value = 33100
$ws.Range("L51").Value = 33100
$ws.Range("N51").Value = -34612
$ws.Range("H74").Value = 19664.36
$ws.Range("I74").Value = 1155.75
$ws.Range("K74").Value = 1155.75
$ws.Range("M74").Value = -281.75
$ws.Range("H77").Value = 19664.36
$ws.Range("I77").Value = 1155.75
$ws.Range("K77").Value = 5778.75
$ws.Range("M77").Value = -1410.75
$ws.Range("H114").Value = 86500
$ws.Range("J114").Value = 86500
$ws.Range("L114").Value = 86500
$ws.Range("N114").Value = -95178
$ws.Range("H132").Value = 1776.5303
$ws.Range("I132").Value = 1421.6885
$ws.Range("K132").Value = 4265.0655
$ws.Range("M132").Value = -1735.0655

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 3182.8
$ws.Range("I25").Value = 1304.6666
$ws.Range("J25").Value = 6000
$ws.Range("K25").Value = 1304.6666
$ws.Range("L25").Value = 6000
$ws.Range("M25").Value = -1069.6666
$ws.Range("N25").Value = -6470
$ws.Range("H95").Value = 9995
$ws.Range("J95").Value = 9995
$ws.Range("L95").Value = 9995
$ws.Range("N95").Value = -15487
$ws.Range("H99").Value = 7146809
$ws.Range("I99").Value = 9527279
$ws.Range("K99").Value = 9527279
$ws.Range("M99").Value = -9525781
$ws.Range("H105").Value = 10417636
$ws.Range("I105").Value = 12500899
$ws.Range("K105").Value = 12500899
$ws.Range("M105").Value = -12499152
$ws.Range("H134").Value = 3260.4226
$ws.Range("I134").Value = 1694.4529
$ws.Range("K134").Value = 5083.3587
$ws.Range("M134").Value = -2548.3587

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1873.7778
$ws.Range("I16").Value = 1553.7142
$ws.Range("K16").Value = 1553.7142
$ws.Range("M16").Value = -1266.7142
$ws.Range("H31").Value = 23823.771
$ws.Range("I31").Value = 3498.3333
$ws.Range("K31").Value = 3498.3333
$ws.Range("M31").Value = -3203.3333
$ws.Range("H34").Value = 23823.771
$ws.Range("I34").Value = 3498.3333
$ws.Range("K34").Value = 3498.3333
$ws.Range("M34").Value = -3296.3333
$ws.Range("H58").Value = 4685.9614
$ws.Range("I58").Value = 5474.4
$ws.Range("K58").Value = 5474.4
$ws.Range("M58").Value = -5271.4
$ws.Range("H86").Value = 7728.3335
$ws.Range("I86").Value = 4247.6665
$ws.Range("J86").Value = 12949.333
$ws.Range("K86").Value = 4247.6665
$ws.Range("L86").Value = 12949.333
$ws.Range("M86").Value = -3124.6665
$ws.Range("N86").Value = -15195.333
$ws.Range("H89").Value = 7728.3335
$ws.Range("I89").Value = 4247.6665
$ws.Range("J89").Value = 12949.333
$ws.Range("K89").Value = 21238.3325
$ws.Range("L89").Value = 64746.665
$ws.Range("M89").Value = -15622.3325
$ws.Range("N89").Value = -75978.66500000001
$ws.Range("H107").Value = 2780.9375
$ws.Range("I107").Value = 2415.4167
$ws.Range("J107").Value = 3877.5
$ws.Range("K107").Value = 2415.4167
$ws.Range("L107").Value = 3877.5
$ws.Range("M107").Value = -495.4167000000002
$ws.Range("N107").Value = -7717.5
$ws.Range("H113").Value = 1873.7778
$ws.Range("I113").Value = 1553.7142
$ws.Range("K113").Value = 1553.7142
$ws.Range("M113").Value = 616.2858000000001
$ws.Range("H122").Value = 2044.6207
$ws.Range("I122").Value = 1963.76
$ws.Range("K122").Value = 5891.28
$ws.Range("M122").Value = -3441.28
$ws.Range("H131").Value = 63899.5
$ws.Range("J131").Value = 63899.5
$ws.Range("L131").Value = 63899.5
$ws.Range("N131").Value = -73979.5
$ws.Range("H132").Value = 47067.605
$ws.Range("I132").Value = 32030.146
$ws.Range("K132").Value = 96090.43799999999
$ws.Range("M132").Value = -93560.43799999999
$ws.Range("H136").Value = 4685.9614
$ws.Range("I136").Value = 5474.4
$ws.Range("K136").Value = 16423.2
$ws.Range("M136").Value = -13873.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20217960
$ws.Range("I4").Value = 23095142
$ws.Range("K4").Value = 69285426
$ws.Range("M4").Value = -69285314
$ws.Range("H92").Value = 706.5714
$ws.Range("I92").Value = 194.14285
$ws.Range("K92").Value = 582.4285500000001
$ws.Range("M92").Value = 665.5714499999999
$ws.Range("H98").Value = 1562.6666
$ws.Range("I98").Value = 1094
$ws.Range("J98").Value = 2500
$ws.Range("K98").Value = 3282
$ws.Range("L98").Value = 7500
$ws.Range("M98").Value = -1784
$ws.Range("N98").Value = -10496
$ws.Range("H113").Value = 3061.257
$ws.Range("I113").Value = 10329
$ws.Range("K113").Value = 30987
$ws.Range("M113").Value = -28817
$ws.Range("H114").Value = 1538.1538
$ws.Range("J114").Value = 1664
$ws.Range("L114").Value = 4992
$ws.Range("N114").Value = -11500
$ws.Range("H117").Value = 3539.5
$ws.Range("I117").Value = 1999
$ws.Range("J117").Value = 3710.6667
$ws.Range("K117").Value = 5997
$ws.Range("L117").Value = 11132.0001
$ws.Range("M117").Value = -2555
$ws.Range("N117").Value = -18016.0001
$ws.Range("H129").Value = 1761.4286
$ws.Range("I129").Value = 1232.75
$ws.Range("J129").Value = 2466.3333
$ws.Range("K129").Value = 3698.25
$ws.Range("L129").Value = 7398.999899999999
$ws.Range("M129").Value = 1301.75
$ws.Range("N129").Value = -17398.9999
$ws.Range("H131").Value = 11265068
$ws.Range("J131").Value = 11908919
$ws.Range("L131").Value = 35726757
$ws.Range("N131").Value = -35736837
$ws.Range("H132").Value = 1389.7778
$ws.Range("I132").Value = 1416.3334
$ws.Range("J132").Value = 1336.6666
$ws.Range("K132").Value = 12747.0006
$ws.Range("L132").Value = 12029.9994
$ws.Range("M132").Value = -10217.0006
$ws.Range("N132").Value = -17089.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 8999
$ws.Range("J4").Value = 8999
$ws.Range("L4").Value = 8999
$ws.Range("N4").Value = -9223
$ws.Range("H39").Value = 68900
$ws.Range("J39").Value = 68900
$ws.Range("L39").Value = 68900
$ws.Range("N39").Value = -69964
$ws.Range("H55").Value = 5875
$ws.Range("I55").Value = 3833.3333
$ws.Range("J55").Value = 12000
$ws.Range("K55").Value = 3833.3333
$ws.Range("L55").Value = 12000
$ws.Range("M55").Value = -3506.3333
$ws.Range("N55").Value = -12654
$ws.Range("H99").Value = 8698.272000000001
$ws.Range("I99").Value = 6073.5557
$ws.Range("K99").Value = 6073.5557
$ws.Range("M99").Value = -3827.5557
$ws.Range("H102").Value = 3486992
$ws.Range("I102").Value = 4116155.2
$ws.Range("K102").Value = 4116155.2
$ws.Range("M102").Value = -4114533.2
$ws.Range("H122").Value = 319947.3
$ws.Range("I122").Value = 469978.84
$ws.Range("K122").Value = 1409936.52
$ws.Range("M122").Value = -1407486.52

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 41437.547
$ws.Range("J22").Value = 1164.1818
$ws.Range("L22").Value = 1164.1818
$ws.Range("N22").Value = -1754.1818
$ws.Range("H27").Value = 41437.547
$ws.Range("J27").Value = 1164.1818
$ws.Range("L27").Value = 1164.1818
$ws.Range("N27").Value = -1378.1818
$ws.Range("H40").Value = 5158.85
$ws.Range("I40").Value = 3148.4167
$ws.Range("J40").Value = 8174.5
$ws.Range("K40").Value = 3148.4167
$ws.Range("L40").Value = 8174.5
$ws.Range("M40").Value = -3012.4167
$ws.Range("N40").Value = -8446.5
$ws.Range("H48").Value = 35000
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H61").Value = 8549489
$ws.Range("I61").Value = 9260613
$ws.Range("K61").Value = 9260613
$ws.Range("M61").Value = -9260411
$ws.Range("H82").Value = 211114510
$ws.Range("I82").Value = 211114510
$ws.Range("K82").Value = 211114510
$ws.Range("M82").Value = -211114149
$ws.Range("H85").Value = 211114510
$ws.Range("I85").Value = 211114510
$ws.Range("K85").Value = 211114510
$ws.Range("M85").Value = -211113262
$ws.Range("H93").Value = 47636750
$ws.Range("I93").Value = 83336090
$ws.Range("J93").Value = 37635.668
$ws.Range("K93").Value = 83336090
$ws.Range("L93").Value = 37635.668
$ws.Range("M93").Value = -83334842
$ws.Range("N93").Value = -40131.668
$ws.Range("H113").Value = 8549489
$ws.Range("I113").Value = 9260613
$ws.Range("K113").Value = 9260613
$ws.Range("M113").Value = -9258443
$ws.Range("H122").Value = 4916.6787
$ws.Range("I122").Value = 2300.7
$ws.Range("J122").Value = 6370
$ws.Range("K122").Value = 6902.099999999999
$ws.Range("L122").Value = 19110
$ws.Range("M122").Value = -4452.099999999999
$ws.Range("N122").Value = -24010
$ws.Range("H132").Value = 6793.8613
$ws.Range("I132").Value = 7318.5864
$ws.Range("K132").Value = 21955.7592
$ws.Range("M132").Value = -19425.7592
$ws.Range("H136").Value = 23283.453
$ws.Range("I136").Value = 38748.715
$ws.Range("J136").Value = 5962.36
$ws.Range("K136").Value = 116246.145
$ws.Range("L136").Value = 17887.08
$ws.Range("M136").Value = -113696.145
$ws.Range("N136").Value = -22987.08

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 15836.167
$ws.Range("I51").Value = 9003.6
$ws.Range("K51").Value = 9003.6
$ws.Range("M51").Value = -8493.6
$ws.Range("H52").Value = 5800
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H54").Value = 1766665.6
$ws.Range("I54").Value = 2525000
$ws.Range("J54").Value = 249997
$ws.Range("K54").Value = 2525000
$ws.Range("L54").Value = 249997
$ws.Range("M54").Value = -2524480
$ws.Range("N54").Value = -251037
$ws.Range("H126").Value = 1513.2812
$ws.Range("I126").Value = 1272.5238
$ws.Range("J126").Value = 1972.909
$ws.Range("K126").Value = 3817.5714
$ws.Range("L126").Value = 5918.727000000001
$ws.Range("M126").Value = -1347.5714
$ws.Range("N126").Value = -10858.727
$ws.Range("H132").Value = 13304587
$ws.Range("I132").Value = 15388285
$ws.Range("J132").Value = 991830
$ws.Range("K132").Value = 46164855
$ws.Range("L132").Value = 2975490
$ws.Range("M132").Value = -46162325
$ws.Range("N132").Value = -2980550
